$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 678
$ws.Range("I33").Value = 666.8333
$ws.Range("J33").Value = 703.125
$ws.Range("K33").Value = 666.8333
$ws.Range("L33").Value = 703.125
$ws.Range("M33").Value = -437.8333
$ws.Range("N33").Value = -1161.125

$ws.Range("H98").Value = 848.3333
$ws.Range("I98").Value = 848.3333
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 848.3333
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 649.6667
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 848.3333
$ws.Range("I122").Value = 848.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2544.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -94.9998999999998
$ws.Range("N122").ClearContents()

$ws.Range("H125").Value = 83333980
$ws.Range("I125").Value = 451
$ws.Range("J125").Value = 125000750
$ws.Range("K125").Value = 4059
$ws.Range("L125").Value = 1125006750
$ws.Range("M125").Value = -1599
$ws.Range("N125").Value = -1125011670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 364416.97
$ws.Range("I61").Value = 291098.88
$ws.Range("J61").Value = 503721.34
$ws.Range("K61").Value = 291098.88
$ws.Range("L61").Value = 503721.34
$ws.Range("M61").Value = -290886.88
$ws.Range("N61").Value = -504145.34

$ws.Range("H74").Value = 152156.77
$ws.Range("I74").Value = 162128.1
$ws.Range("K74").Value = 162128.1
$ws.Range("M74").Value = -161254.1

$ws.Range("H77").Value = 152156.77
$ws.Range("I77").Value = 162128.1
$ws.Range("K77").Value = 810640.5
$ws.Range("M77").Value = -806272.5

$ws.Range("H110").Value = 1334.8572
$ws.Range("I110").Value = 1556.8
$ws.Range("J110").Value = 780
$ws.Range("K110").Value = 1556.8
$ws.Range("L110").Value = 780
$ws.Range("M110").Value = 488.2
$ws.Range("N110").Value = -4870

$ws.Range("H122").Value = 4137.1953
$ws.Range("I122").Value = 4220.2188
$ws.Range("K122").Value = 12660.6564
$ws.Range("M122").Value = -10210.6564

$ws.Range("H136").Value = 364416.97
$ws.Range("I136").Value = 291098.88
$ws.Range("J136").Value = 503721.34
$ws.Range("K136").Value = 873296.64
$ws.Range("L136").Value = 1511164.02
$ws.Range("M136").Value = -870746.64
$ws.Range("N136").Value = -1516264.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1129.75
$ws.Range("I94").Value = 883.3333
$ws.Range("J94").Value = 1277.6
$ws.Range("K94").Value = 883.3333
$ws.Range("L94").Value = 1277.6
$ws.Range("M94").Value = -432.3333
$ws.Range("N94").Value = -2179.6

$ws.Range("H134").Value = 2724.3242
$ws.Range("I134").Value = 2244.9583
$ws.Range("J134").Value = 3609.3076
$ws.Range("K134").Value = 6734.874899999999
$ws.Range("L134").Value = 10827.9228
$ws.Range("M134").Value = -4199.874899999999
$ws.Range("N134").Value = -15897.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 102561.2
$ws.Range("I99").Value = 168502
$ws.Range("J99").Value = 3650
$ws.Range("K99").Value = 168502
$ws.Range("L99").Value = 3650
$ws.Range("M99").Value = -167004
$ws.Range("N99").Value = -6646

$ws.Range("H122").Value = 1547.6471
$ws.Range("I122").Value = 838.875
$ws.Range("J122").Value = 2177.6667
$ws.Range("K122").Value = 2516.625
$ws.Range("L122").Value = 6533.000100000001
$ws.Range("M122").Value = -66.625
$ws.Range("N122").Value = -11433.0001

$ws.Range("H126").Value = 102561.2
$ws.Range("I126").Value = 168502
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 505506
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -503036
$ws.Range("N126").Value = -15890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H52").Value = 10000
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 10000
$ws.Range("N52").Value = -10518

$ws.Range("H102").Value = 4406.905
$ws.Range("I102").Value = 2174.4119
$ws.Range("J102").Value = 13895
$ws.Range("K102").Value = 2174.4119
$ws.Range("L102").Value = 13895
$ws.Range("M102").Value = -552.4119000000001
$ws.Range("N102").Value = -17139

$ws.Range("H122").Value = 1172.6364
$ws.Range("I122").Value = 1211.1111
$ws.Range("J122").Value = 999.5
$ws.Range("K122").Value = 3633.3333
$ws.Range("L122").Value = 2998.5
$ws.Range("M122").Value = -1183.3333
$ws.Range("N122").Value = -7898.5

$ws.Range("H126").Value = 2014.6774
$ws.Range("I126").Value = 1758.1305
$ws.Range("J126").Value = 2752.25
$ws.Range("K126").Value = 5274.3915
$ws.Range("L126").Value = 8256.75
$ws.Range("M126").Value = -2804.3915
$ws.Range("N126").Value = -13196.75

$ws.Range("H132").Value = 4335.147
$ws.Range("I132").Value = 4944.6313
$ws.Range("K132").Value = 14833.8939
$ws.Range("M132").Value = -12303.8939

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2280.1428
$ws.Range("I7").Value = 1830.6
$ws.Range("J7").Value = 3404
$ws.Range("K7").Value = 1830.6
$ws.Range("L7").Value = 3404
$ws.Range("M7").Value = -1718.6
$ws.Range("N7").Value = -3628

$ws.Range("H122").Value = 2841.8333
$ws.Range("I122").Value = 2257.4285
$ws.Range("J122").Value = 3660
$ws.Range("K122").Value = 6772.2855
$ws.Range("L122").Value = 10980
$ws.Range("M122").Value = -4322.2855
$ws.Range("N122").Value = -15880

$ws.Range("H126").Value = 2280.1428
$ws.Range("I126").Value = 1830.6
$ws.Range("J126").Value = 3404
$ws.Range("K126").Value = 5491.799999999999
$ws.Range("L126").Value = 10212
$ws.Range("M126").Value = -3021.799999999999
$ws.Range("N126").Value = -15152

$ws.Range("H136").Value = 2502.9404
$ws.Range("I136").Value = 1249.0312
$ws.Range("J136").Value = 6515.45
$ws.Range("K136").Value = 3747.0936
$ws.Range("L136").Value = 19546.35
$ws.Range("M136").Value = -1197.0936
$ws.Range("N136").Value = -24646.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14950
$ws.Range("J54").Value = 14950
$ws.Range("L54").Value = 14950
$ws.Range("N54").Value = -15990

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H81").Value = 2242.1052
$ws.Range("I81").Value = 1492.3077
$ws.Range("J81").Value = 3866.6667
$ws.Range("K81").Value = 2984.6154
$ws.Range("L81").Value = 7733.3334
$ws.Range("M81").Value = -1923.6154
$ws.Range("N81").Value = -9855.3334

$ws.Range("H84").Value = 2242.1052
$ws.Range("I84").Value = 1492.3077
$ws.Range("J84").Value = 3866.6667
$ws.Range("K84").Value = 14923.077
$ws.Range("L84").Value = 38666.667
$ws.Range("M84").Value = -9619.077000000001
$ws.Range("N84").Value = -49274.667

$ws.Range("H126").Value = 2100.147
$ws.Range("I126").Value = 2059.9048
$ws.Range("J126").Value = 2165.1538
$ws.Range("K126").Value = 6179.714399999999
$ws.Range("L126").Value = 6495.4614
$ws.Range("M126").Value = -3709.714399999999
$ws.Range("N126").Value = -11435.4614

$ws.Range("H136").Value = 12799036
$ws.Range("I136").Value = 16146297
$ws.Range("J136").Value = 591376.75
$ws.Range("K136").Value = 48438891
$ws.Range("L136").Value = 1774130.25
$ws.Range("M136").Value = -48436341
$ws.Range("N136").Value = -1779230.25
